$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$names = @(
  '126 Racecourse Road Public Housing Tower Flemington',
  '139 Highett St Apartment Complex Richmond',
  '3175 The Bays Aged Care Facility Hastings',
  '3535 Opal Meadow Heights Aged Care Community Meadow Heights',
  'Al Haj Halal Meats Glenroy',
  'Al-Taqwa College Truganina',
  'Allbright Manor Aged Care Croydon North Tier 1A',
  'Australia Post Distribution Centre Sunshine West',
  'Australian Lamb Colac East',
  'Baker Bleu Caulfield North',
  'Baxter Foods Australia Campbellfield',
  'CFMEU Melbourne Office',
  'CS Square Caroline Springs',
  'Cafe Roco Dandenong',
  'Campbellfield Ford Complex Vaccination Clinic Campbellfield',
  'Cardinia Lakes Early Learning Centre Pakenham',
  'Caroline Springs Police Station',
  'Carton Finishing Pty. Ltd. Campbellfield',
  'Chemist Warehouse Campbellfield DC',
  'Chemist Warehouse Fillo Drive Somerton',
  'City of Wyndham Community',
  'Coles Campbellfield Plaza Campbellfield',
  'Coles Coburg North Village',
  'Coles Pakenham Place Shopping Centre',
  'Coles Roxburgh Village Roxburgh Park',
  'Community Kids Bayswater Early Education Centre Bayswater North',
  'Construction Site 1 Warde Street Footscray',
  'Construction Site Olea Apartment Caulfield North',
  'Costco Wholesale Epping',
  'Crusader Caravans Epping',
  'Dandenong Police Station Dandenong',
  'DayHab Rehabilitation Treatment Centre Ringwood East',
  'Direct Freight Express Campbellfield',
  'Disability Residence Life without Barriers Ashwood',
  'Don Watson Coldstore Derrimut',
  'Epworth Healthcare Epworth Richmond Emergency Department',
  'Ermha365 Residential Disability Care Services Paperbark St Doveton',
  'FedEx Station Melbourne Airport',
  'Fine Food Holdings Pty Ltd Dandenong South',
  'Fonterra Manufacturing Workplace Campbellfield',
  'General Foods Campbellfield',
  'Gladstone Parade Early Learning & Kinder Glenroy',
  'Goodstart Early Learning Altona',
  'Green Leaves Early Learning Cairnlea',
  'Green Leaves Early Learning Centre Highlands Craigieburn',
  'Greenvale Primary School',
  'HEI Schools Emerald Early Learning Centre Emerald',
  'Hamilton Marino 236 Jasper Road McKinnon',
  'Hello Fresh Warehouse Ravenhall',
  'ISS Factory Level 1 Terminal 2 Melbourne Airport Tullamarine',
  'Ibis Kingsgate Hotel Melbourne',
  'Industrial Galvanizers Valmont Coatings Campbellfield',
  'Inghams Enterprises Thomastown',
  'Kippers Seafood Werribee',
  'Kool Kidz Childcare Narre Warren',
  'Lantmannen Unibake Australia Mordialloc',
  'Linfox Somerton National Distribution Centre Somerton',
  'Mecca Distribution Centre Warehouse Melbourne Airport',
  'Melbourne Assessment Prison West Melbourne',
  'Melbourne Metropolitan Remand Centre Ravenhall',
  'Melbourne West Police Station Docklands',
  'Mill Park Police Station Mill Park',
  'MyCentre Childcare Broadmeadows',
  'National Gallery of Victoria Melbourne',
  'Nido Early School Ascot Vale',
  'Nido Early School Glenroy',
  'Northern Health Northern Hospital Epping Emergency Department Tier 1B',
  'Northern Health The Northern Hospital Epping',
  'OnQ Plumbing and Excavations Craigieburn',
  'Oporto Coolaroo',
  'Oscar Romero Catholic Primary School Craigieburn',
  'Our Lady Help of Christian''s Primary School Brunswick East',
  'Pacific Meat Thomastown',
  'Private Residence Daycare Allumba Way Wollert',
  'Ravenhall Correctional Centre Ravenhall',
  'Richmond Quarter 261-271 Bridge Road Construction Site Richmond',
  'Sacca''s Fruit World Broadmeadows Central Shopping Centre',
  'Sharpline Stainless Steel Coburg North',
  'St Margaret''s Primary School OSHC Maribyrnong',
  'St Vincents Hospital Emergency Department Melbourne',
  'Tek Foods Somerton',
  'The Huntly-Goornong Rail Works',
  'The Royal Children''s Hospital Melbourne Emergency Department Parkville Tier 1B',
  'The Royal Melbourne Hospital Parkville',
  'The Royal Melbourne Hospital Parkville Emergency Department',
  'The Royal Melbourne Hospital Ward 6SE Parkville',
  'The Royal Talbot Rehabilitation Centre Kew',
  'ThorwestenCabinets Pakenham',
  'Truganina Early Learning Centre Truganina',
  'Visy Recycling Springvale',
  'Wallaby Childcare Wollert',
  'Werribee Mercy Hospital Emergency Department',
  'Western Health Footscray Hospital Emergency Department',
  'Western Health Sunshine Hospital Emergency Department',
  'Western Health Sunshine Hospital GEM Ward St Albans',
  'Yara Childcare Centre Truganina'
)

$counts = @(
  6,
  9,
  6,
  18,
  11,
  8,
  5,
  5,
  9,
  7,
  5,
  5,
  12,
  6,
  11,
  5,
  8,
  12,
  5,
  9,
  6,
  7,
  18,
  6,
  5,
  17,
  5,
  11,
  25,
  22,
  9,
  7,
  5,
  5,
  5,
  7,
  9,
  15,
  10,
  9,
  14,
  7,
  9,
  5,
  16,
  5,
  5,
  9,
  7,
  9,
  6,
  14,
  6,
  5,
  15,
  6,
  9,
  9,
  9,
  7,
  6,
  8,
  12,
  7,
  12,
  23,
  50,
  19,
  10,
  6,
  5,
  9,
  5,
  8,
  10,
  11,
  6,
  6,
  12,
  14,
  5,
  6,
  17,
  6,
  6,
  14,
  12,
  11,
  7,
  23,
  18,
  15,
  8,
  14,
  6,
  5
)

for ($i = 0; $i -lt $names.Count; $i++) {
  $row = $i + 2
  $ws.Cells.Item($row, 1).Value = $names[$i]
  $ws.Cells.Item($row, 2).Value = $counts[$i]
}
